$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that was bumped by one day
# (45203 -> 45204) for every data row (rows 2 through 526).
$ws.Range("C2:C526").Value = 45204
